$wb = $excel.ActiveWorkbook

# --- Insert new worksheet "sleepdiary2" right after "Sleep Diary" ---
$sheetDiary = $wb.Worksheets.Item("Sleep Diary")
$newSheet = $wb.Worksheets.Add($null, $sheetDiary)
$newSheet.Name = "sleepdiary2"
$ws = $wb.Worksheets.Item("sleepdiary2")

# --- Header row ---
$ws.Range("A1").Value = "Dato"
# Copy formatting (fill/border) from the "Sleep Diary" date header cell
$sheetDiary.Range("A3").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("B1").Value = "q1"
$ws.Range("C1").Value = "q2"
$ws.Range("D1").Value = "q3"
$ws.Range("E1").Value = "q4"
$ws.Range("F1").Value = "q5"
$ws.Range("G1").Value = "q6"
$ws.Range("H1").Value = "q7"
$ws.Range("I1").Value = "q8"

# --- Date column (A2:A15) ---
$dates = @(44263,44264,44265,44266,44267,44268,44269,44270,44271,44272,44273,44274,44275,44276)
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}
# Copy date-cell formatting (fill/border/number format) from "Sleep Diary"
$sheetDiary.Range("A4").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)

# --- Time-of-day values (columns B, C, G, H) ---
$colB = @(0.89583333333333337,0.91666666666666663,0.9375,0.89583333333333337,0.91666666666666663,0.94791666666666663,0.020833333333333332,0.91666666666666663,0.91666666666666663,0.9375,0.91666666666666663,0.020833333333333332,0.083333333333333329,0.083333333333333329)
$colC = @(0.95833333333333337,0.95833333333333337,0.97916666666666663,0.91666666666666663,0.95833333333333337,0.97916666666666663,0.020833333333333332,0.9375,0.9375,0.95833333333333337,0.95833333333333337,0.083333333333333329,0.10416666666666667,0.10416666666666667)
$colG = @(0.33333333333333331,0.3125,0.33333333333333331,0.33333333333333331,0.27083333333333331,0.35416666666666669,0.32291666666666669,0.32291666666666669,0.33055555555555555,0.3298611111111111,0.29166666666666669,0.33333333333333331,0.36805555555555558,0.2951388888888889)
$colH = @(0.34375,0.33333333333333331,0.33333333333333331,0.33333333333333331,0.28125,0.375,0.32291666666666669,0.33333333333333331,0.33680555555555558,0.33333333333333331,0.33333333333333331,0.37361111111111112,0.375,0.33124999999999999)

for ($i = 0; $i -lt $colB.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
    $ws.Cells.Item($row, 7).Value = $colG[$i]
    $ws.Cells.Item($row, 8).Value = $colH[$i]
}
$ws.Range("B2:C15").NumberFormat = "h:mm"
$ws.Range("G2:H15").NumberFormat = "h:mm"

# --- Plain numeric columns D, E, F, I ---
$ws.Cells.Item(8, 4).Value = 30
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(13, 4).Value = 5
$ws.Cells.Item(14, 4).Value = 2

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(9, 5).Value = 2

$ws.Cells.Item(8, 6).Value = 30
$ws.Cells.Item(9, 6).Value = 5

$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(3, 9).Value = 3
$ws.Cells.Item(4, 9).Value = 4
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(7, 9).Value = 4
$ws.Cells.Item(8, 9).Value = 2
$ws.Cells.Item(9, 9).Value = 5
$ws.Cells.Item(10, 9).Value = 5
$ws.Cells.Item(11, 9).Value = 4
$ws.Cells.Item(12, 9).Value = 3
$ws.Cells.Item(13, 9).Value = 4
$ws.Cells.Item(14, 9).Value = 4
$ws.Cells.Item(15, 9).Value = 3

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 10.71

# --- Selection / view state for the new sheet ---
$ws.Range("G15").Select()

# --- Adjust selection/view on the other sheets that changed ---
$wsQuiz = $wb.Worksheets.Item("Sleep Quiz")
$wsQuiz.Select()
$wsQuiz.Range("C9").Select()

$sheetDiary.Select()
$sheetDiary.Range("A3:A17").Select()
$sheetDiary.Range("A17").Activate()

$ws.Select()
